$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.017234
$ws.Range("H2").Value = 0.051702
$ws.Range("I2").Value = 0.001089091024602478
$ws.Range("J2").Value = 0.001089091024602478
$ws.Range("M2").Value = 560.2199806666666
$ws.Range("N2").Value = 1680.659942
$ws.Range("O2").Value = 0.6936344353529325
$ws.Range("P2").Value = 0.6936344353529326
$ws.Range("Q2").Value = 9.654831146809332
$ws.Range("R2").Value = 86.893480321284
$ws.Range("S2").Value = 0.0007554310378980867
$ws.Range("T2").Value = 0.0007554310378980868

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.017234
$ws.Range("H3").Value = 0.051702
$ws.Range("I3").Value = 0.001089091024602478
$ws.Range("J3").Value = 0.001089091024602478
$ws.Range("O3").Value = 0.06994956469466522
$ws.Range("P3").Value = 0.06994956469466522
$ws.Range("Q3").Value = 0.9736414478559999
$ws.Range("R3").Value = 8.762773030704
$ws.Range("S3").Value = 0.000076181443083810277783034237
$ws.Range("T3").Value = 0.000076181443083810277783034237

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.017234
$ws.Range("H4").Value = 0.051702
$ws.Range("I4").Value = 0.001089091024602478
$ws.Range("J4").Value = 0.001089091024602478
$ws.Range("O4").Value = 0.2364159999524024
$ws.Range("P4").Value = 0.2364159999524024
$ws.Range("Q4").Value = 3.290719784958667
$ws.Range("R4").Value = 29.616478064628
$ws.Range("S4").Value = 0.0002574785436205813
$ws.Range("T4").Value = 0.0002574785436205813

# Row 5
$ws.Range("I5").Value = 0.991605223525074
$ws.Range("J5").Value = 0.991605223525074
$ws.Range("M5").Value = 560.2199806666666
$ws.Range("N5").Value = 1680.659942
$ws.Range("O5").Value = 0.6936344353529325
$ws.Range("P5").Value = 0.6936344353529326
$ws.Range("Q5").Value = 8790.616010193617
$ws.Range("R5").Value = 79115.54409174256
$ws.Range("S5").Value = 0.6878115293128331
$ws.Range("T5").Value = 0.6878115293128332

# Row 6
$ws.Range("I6").Value = 0.991605223525074
$ws.Range("J6").Value = 0.991605223525074
$ws.Range("O6").Value = 0.06994956469466522
$ws.Range("P6").Value = 0.06994956469466522
$ws.Range("Q6").Value = 886.4896723273654
$ws.Range("R6").Value = 7978.407050946289
$ws.Range("S6").Value = 0.06936235373453513
$ws.Range("T6").Value = 0.06936235373453513

# Row 7
$ws.Range("I7").Value = 0.991605223525074
$ws.Range("J7").Value = 0.991605223525074
$ws.Range("O7").Value = 0.2364159999524024
$ws.Range("P7").Value = 0.2364159999524024
$ws.Range("Q7").Value = 2996.163639410547
$ws.Range("R7").Value = 26965.47275469492
$ws.Range("S7").Value = 0.2344313404777058
$ws.Range("T7").Value = 0.2344313404777058

# Row 8
$ws.Range("G8").Value = 0.1156066666666667
$ws.Range("I8").Value = 0.007305685450323614
$ws.Range("J8").Value = 0.007305685450323614
$ws.Range("M8").Value = 560.2199806666666
$ws.Range("N8").Value = 1680.659942
$ws.Range("O8").Value = 0.6936344353529325
$ws.Range("P8").Value = 0.6936344353529326
$ws.Range("Q8").Value = 64.76516456493778
$ws.Range("R8").Value = 582.88648108444
$ws.Range("S8").Value = 0.005067475002201355
$ws.Range("T8").Value = 0.005067475002201355

# Row 9
$ws.Range("G9").Value = 0.1156066666666667
$ws.Range("I9").Value = 0.007305685450323614
$ws.Range("J9").Value = 0.007305685450323614
$ws.Range("O9").Value = 0.06994956469466522
$ws.Range("P9").Value = 0.06994956469466522
$ws.Range("Q9").Value = 6.531243026293334
$ws.Range("R9").Value = 58.78118723664
$ws.Range("S9").Value = 0.0005110295170462861
$ws.Range("T9").Value = 0.000511029517046286

# Row 10
$ws.Range("G10").Value = 0.1156066666666667
$ws.Range("I10").Value = 0.007305685450323614
$ws.Range("J10").Value = 0.007305685450323614
$ws.Range("O10").Value = 0.2364159999524024
$ws.Range("P10").Value = 0.2364159999524024
$ws.Range("S10").Value = 0.001727180931075974
$ws.Range("T10").Value = 0.001727180931075974
